# Remove the gimme_boxes frame type; move this boolean (and the
# succeedsTransport argument) into the hello frame.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Protocol")
$ws.Activate()

# 1. Update the "hello, helloData" comment cell (F5) to document the new
#    "g" (succeedsTransport) key that now lives in the hello frame, taking
#    over the duty previously served by the standalone gimme_boxes frame.
$helloComment = "helloData: {""n"": transportNumber, ""v: protocolVersion, ""t"": httpFormat, ""w"": requestNewStream, ""i"": streamId, ""c"": credentialsData, ""p"": needPaddingBytes, ""r"": maxReceiveBytes, ""m"": maxOpenTime, ""o"": readOnlyOnce: ""a"": useMyTcpAcks, ""g"": succeedsTransport}. Presence of ""g"" means ""give me boxes, server"". If succeedsTransport != null, temporarily assume that all boxes written to #<succeedsTransport> were SACKed." + [char]10 + "Only C2S because hello is used by the client to identify itself to the server, and set critical transport parameters. XXX TODO: perhaps ackMode: 0 - require Minerva-level SACKs, 1 - use my TCP acks, 2 - assume everything written is received"
$ws.Range("F5").Value = $helloComment

# 2. Remove the gimme_boxes row (row 7) entirely: the gimme_boxes frame type
#    no longer exists as its own row in the table, it's folded into hello.
$ws.Range("A7:F7").Clear()
$ws.Rows.Item(7).RowHeight = 12.75

# 3. Reflect the row-7 deletion in the sheet's current selection / active
#    cell (the row that used to hold gimme_boxes is now selected/empty).
$ws.Range("A7:F7").Select()
